# Mise a jour de certains champs de Modules et de Professeurs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Enseignant"/"Nombre d'heures" header text into
# "Composants" (col D) / "Chef  Module" (col C).
# D1 is assigned first so the shared-string table keeps the same
# index order as the target file (Composants at index 2, Chef Module at index 3).
$ws.Range("D1").Value = "Composants"
$ws.Range("C1").Value = "Chef  Module"

# Column widths for C and D (values chosen so the engine's rounding
# lands on the target stored widths of 35 and ~24.57 characters).
$ws.Columns.Item(3).ColumnWidth = 34.166666666666664
$ws.Columns.Item(4).ColumnWidth = 23.666666666666668

# Update the active selection to E8.
$ws.Range("E8").Select()
